$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6222347617149353
$ws.Range("B1").Value = 0.8970164060592651
$ws.Range("C1").Value = 5.984857559204102
$ws.Range("D1").Value = 3.385230302810669
$ws.Range("E1").Value = 0.9762418270111084
